$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$students = @(
    ,@(10, "LuisaRamirez")
    ,@(11, "JessicaArcila")
    ,@(12, "ManuelaSotoMarin")
    ,@(13, "JoselynArangoBedoya")
    ,@(14, "ManuelNorenaGuevera")
    ,@(15, "DannaMolinaZapata")
    ,@(16, "YuriJohanaVinascoFonseca")
    ,@(17, "LauraAndreaTrujilloRestrepo")
    ,@(18, "JulianaRioMartinez")
    ,@(19, "ValentinaRamirezCuesta")
    ,@(20, "HaroldValenciaGonzales")
    ,@(21, "GeraldinPalacioMunos")
    ,@(22, "MichellyHurtadoJimenez")
    ,@(23, "DamaraGiraldoBolivar")
    ,@(24, "SaraAcevedoMarin")
    ,@(25, "YuranyRojasPulgarin")
    ,@(26, "JulianaQuinteroArboleda")
    ,@(27, "DuvanEsneiderDiazMontoya")
    ,@(28, "YeilyFunesOsorio")
    ,@(29, "CrisLauraPadilla")
    ,@(30, "MelanyMorales")
    ,@(31, "SofiaVeraMartina")
    ,@(32, "MarianaValenciaZapata")
    ,@(33, "JuanCamiloVasqeuz")
    ,@(34, "SebastianMejiaVasquez")
    ,@(35, "SantiagoMejiaAcevedo")
    ,@(36, "EadySantiagoLondono")
    ,@(37, "SalomeSisquiarcoRios")
    ,@(38, "JhonatanMontanaSolano")
    ,@(39, "SantiagoParraOsorio")
    ,@(40, "YenniferJimenezAguila")
    ,@(41, "SamuelYepesOsorio")
)

$startRow = 12
$styleSource = $ws.Cells.Item(11, 1)
$styleSource.Copy()

for ($i = 0; $i -lt $students.Count; $i++) {
    $r = $startRow + $i
    $num = $students[$i][0]
    $name = $students[$i][1]
    $target = $ws.Cells.Item($r, 1)
    $target.PasteSpecial(-4122)
    $target.Value = $num
    $ws.Cells.Item($r, 2).Value = $name
    $ws.Cells.Item($r, 3).Value = "No vino"
}
